$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $c = $ws.Range($addr)
    $c.Value = "'" + $val
    $c.Style = "Normal"
}

Set-TextCell "D2" '26.181.30'
$ws.Range("E2").Value = '  -0.58%  '

Set-TextCell "D3" '1.588.26'
$ws.Range("E3").Value = '  +0.00%  '

$ws.Range("E4").Value = '  -0.09%  '

Set-TextCell "D5" '211.56'
$ws.Range("E5").Value = '  +0.68%  '

Set-TextCell "D6" '0.502'
$ws.Range("E6").Value = '  -0.48%  '

$ws.Range("E7").Value = '  -0.06%  '

Set-TextCell "D8" '0.244'
$ws.Range("E8").Value = '  -0.21%  '

$ws.Range("E9").Value = '  -1.03%  '

Set-TextCell "D10" '19.17'
$ws.Range("E10").Value = '  -1.65%  '

Set-TextCell "D11" '0.0846'
$ws.Range("E11").Value = '  -0.02%  '

Set-TextCell "D12" '1.811.88'
$ws.Range("E12").Value = '  -0.03%  '

Set-TextCell "D13" '1.588.86'
$ws.Range("E13").Value = '  -0.14%  '

Set-TextCell "D14" '4.00'
$ws.Range("E14").Value = '  -1.61%  '

Set-TextCell "D15" '0.512'
$ws.Range("E15").Value = '  -1.28%  '

Set-TextCell "D16" '63.64'
$ws.Range("E16").Value = '  -1.04%  '

Set-TextCell "D17" '26.200.84'
$ws.Range("E17").Value = '  -0.54%  '

Set-TextCell "D18" '0.0₃0724'
$ws.Range("E18").Value = '  -0.32%  '

Set-TextCell "D19" '7.37'
$ws.Range("E19").Value = '  -1.16%  '

Set-TextCell "D20" '213.76'
$ws.Range("E20").Value = '  +1.46%  '

$ws.Range("E21").Value = '  -0.11%  '

Set-TextCell "D22" '4.23'
$ws.Range("E22").Value = '  -0.84%  '

Set-TextCell "D23" '8.96'
$ws.Range("E23").Value = '  +0.41%  '

Set-TextCell "D24" '2.11'
$ws.Range("E24").Value = '  -1.53%  '

Set-TextCell "D25" '144.43'
$ws.Range("E25").Value = '  -0.42%  '

$ws.Range("E26").Value = '  -0.10%  '

Set-TextCell "D27" '6.97'
$ws.Range("E27").Value = '  -1.21%  '

$ws.Range("E28").Value = '  -1.23%  '

Set-TextCell "D29" '15.06'
$ws.Range("E29").Value = '  -1.10%  '

$ws.Range("E30").Value = '  -2.20%  '

$ws.Range("E31").Value = '  +0.29%  '

$ws.Range("E32").Value = '  -1.17%  '

Set-TextCell "D33" '1.409.37'
$ws.Range("E33").Value = '  +8.08%  '

$ws.Range("E34").Value = '  -1.31%  '

$ws.Range("E35").Value = '  -0.61%  '

$ws.Range("B36").Value = 'LidoDAOToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextCell "D36" '1.45'
$ws.Range("E36").Value = '  -1.29%  '

$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell "D37" '0.587'
$ws.Range("E37").Value = '  -3.84%  '

$ws.Range("E38").Value = '  -1.36%  '

Set-TextCell "D39" '5.89'
$ws.Range("E39").Value = '  +4.85%  '

$ws.Range("E40").Value = '  +1.48%  '

$ws.Range("E41").Value = '  -0.06%  '

Set-TextCell "D42" '0.945'
$ws.Range("E42").Value = '  -13.31%  '

$ws.Range("E43").Value = '  -0.37%  '

$ws.Range("E44").Value = '  +0.03%  '

Set-TextCell "D45" '1.723.74'
$ws.Range("E45").Value = '  -0.08%  '

Set-TextCell "D46" '60.91'
$ws.Range("E46").Value = '  -2.51%  '

Set-TextCell "D47" '85.85'
$ws.Range("E47").Value = '  -2.15%  '

$ws.Range("E48").Value = '  -0.16%  '

$ws.Range("E49").Value = '  -0.69%  '

Set-TextCell "D50" '0.0956'
$ws.Range("E50").Value = '  -2.31%  '

Set-TextCell "D51" '0.999'
$ws.Range("E51").Value = '  -0.08%  '
